$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename
$ws.Range("O1").Value = "F1 train"

# Row 2 - O column update
$ws.Range("O2").Value = 0.9166666666666666

# Row 3 - O column update
$ws.Range("O3").Value = 1

# Row 4 - O column update
$ws.Range("O4").Value = 0.8767123287671232

# Row 5 - O column update
$ws.Range("O5").Value = 1

# Row 6 - parameter and multiple metric updates
$ws.Range("C6").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 9
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 0.35
$ws.Range("J6").Value = 0.4347826086956522
$ws.Range("K6").Value = 0.5555555555555556
$ws.Range("L6").Value = 0.3571428571428572
$ws.Range("M6").Value = 0.1818181818181818
$ws.Range("N6").Value = 0.5555555555555556
$ws.Range("O6").Value = 0.631578947368421

# Row 7 - O column update
$ws.Range("O7").Value = 1

# Row 8 - O column update
$ws.Range("O8").Value = 1

# Row 9 - O column update
$ws.Range("O9").Value = 0.9859154929577465

# Row 10 - O column update
$ws.Range("O10").Value = 1

# Row 11 - parameter and multiple metric updates
$ws.Range("C11").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 0.6
$ws.Range("J11").Value = 0.5555555555555556
$ws.Range("K11").Value = 0.5555555555555556
$ws.Range("L11").Value = 0.5555555555555556
$ws.Range("M11").Value = 0.6363636363636364
$ws.Range("N11").Value = 0.5555555555555556
$ws.Range("O11").Value = 0.7536231884057971

# Row 12 - O column update
$ws.Range("O12").Value = 0.9565217391304348

# Row 13 - O column update
$ws.Range("O13").Value = 0.9705882352941176

# Row 14 - O column update
$ws.Range("O14").Value = 1

# Row 15 - O column update
$ws.Range("O15").Value = 1

# Row 16 - parameter and multiple metric updates
$ws.Range("C16").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 0.45
$ws.Range("J16").Value = 0.4761904761904762
$ws.Range("K16").Value = 0.5555555555555556
$ws.Range("L16").Value = 0.4166666666666667
$ws.Range("M16").Value = 0.3636363636363636
$ws.Range("N16").Value = 0.5555555555555556
$ws.Range("O16").Value = 0.7792207792207793
